$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the shared-string "BNT" value in A2 with the numeric value 1
$ws.Range("A2").Value = 1

# Move the active selection from I2 to A2
$ws.Range("A2").Select()
